$wb = $excel.ActiveWorkbook

# --- Sheet: info ---
$ws = $wb.Worksheets.Item("info")
$ws.Range("E2").Value = " 0 // 0 - expertos; 1 - Igual importancia; 2 - Enfoque Ambiental; 3 - Enfoque Económico; 4 - Enfoque Técnico"

# --- Sheet: alternative_info ---
$ws = $wb.Worksheets.Item("alternative_info")
$ws.Range("G4").Value = 640.793947259624
$ws.Range("J4").Value = 223.5449862521488
$ws.Range("L4").Value = 0.0007506139335251635
$ws.Range("M4").Value = 183.2282421272689
$ws.Range("G6").Value = 1067.989912099374
$ws.Range("J6").Value = 278.00170466371
$ws.Range("L6").Value = 0.0006455764403994414
$ws.Range("M6").Value = 175.0574402764045
$ws.Range("G7").Value = 640.793947259624
$ws.Range("J7").Value = 1023.077989380044
$ws.Range("L7").Value = 0.6537430848820482
$ws.Range("M7").Value = 122.2633983069356
$ws.Range("G9").Value = 640.793947259624
$ws.Range("J9").Value = 153.5523892969807
$ws.Range("L9").Value = 0.0005999811623896996
$ws.Range("M9").Value = 171.1756420043254
$ws.Range("G11").Value = 1067.989912099374
$ws.Range("J11").Value = 208.0091077085419
$ws.Range("L11").Value = 0.0004769226130545392
$ws.Range("M11").Value = 161.6726742933793
$ws.Range("G12").Value = 640.793947259624
$ws.Range("J12").Value = 953.0853924248759
$ws.Range("L12").Value = 0.7506442802556406
$ws.Range("M12").Value = 105.8177015515107
$ws.Range("G14").Value = 1067.989912099374
$ws.Range("J14").Value = 1007.542110836437
$ws.Range("L14").Value = 0.7043540549013443
$ws.Range("M14").Value = 107.4352358608177
$ws.Range("G15").Value = 640.793947259624
$ws.Range("J15").Value = 1485.968526752771
$ws.Range("L15").Value = 0.8341918991733863
$ws.Range("M15").Value = 95.49808888516999
$ws.Range("G17").Value = 1067.989912099374
$ws.Range("J17").Value = 138.0165107533738
$ws.Range("L17").Value = 0.000009280258328686273
$ws.Range("M17").Value = 124.5594727710923
$ws.Range("G18").Value = 640.793947259624
$ws.Range("J18").Value = 883.0927954697078
$ws.Range("L18").Value = 0.8813317657511304
$ws.Range("M18").Value = 81.02873272804024
$ws.Range("G20").Value = 1067.989912099374
$ws.Range("J20").Value = 937.549513881269
$ws.Range("L20").Value = 0.8181862245994608
$ws.Range("M20").Value = 85.51782949900341
$ws.Range("G21").Value = 640.793947259624
$ws.Range("J21").Value = 1415.975929797603
$ws.Range("L21").Value = 0.9257343445581757
$ws.Range("M21").Value = 77.82034940911073
$ws.Range("G22").Value = 2135.979824198747
$ws.Range("J22").Value = 274.3666084183291
$ws.Range("G23").Value = 1067.989912099374
$ws.Range("J23").Value = 1470.432648209164
$ws.Range("L23").Value = 0.8827955896432095
$ws.Range("M23").Value = 81.04021835940149

# --- Sheet: alternatives_norm ---
$ws = $wb.Worksheets.Item("alternatives_norm")
$ws.Range("B2").Value = 0.01512274081242303
$ws.Range("D2").Value = 0.004043673032953418
$ws.Range("E2").Value = 0.02295223734114445
$ws.Range("B3").Value = 0.007606110179460675
$ws.Range("D3").Value = 0.004061775569168601
$ws.Range("E3").Value = 0.02301992478183012
$ws.Range("B4").Value = 0.01200344868318235
$ws.Range("D4").Value = 0.005387154237815462
$ws.Range("E4").Value = 0.02542144699929226
$ws.Range("B5").Value = 0.05053775241603135
$ws.Range("D5").Value = 0.000005652575252862016
$ws.Range("E5").Value = 0.03821806946560664
$ws.Range("B6").Value = 0.01492755105678861
$ws.Range("D6").Value = 0.006263662643034885
$ws.Range("E6").Value = 0.02660799243184001
$ws.Range("B7").Value = 0.05493509091975302
$ws.Range("D7").Value = 0.000006185416146593732
$ws.Range("E7").Value = 0.03809747733592668
$ws.Range("B8").Value = 0.07915139061153695
$ws.Range("D8").Value = 0.000005006886017537524
$ws.Range("E8").Value = 0.04339396233816271
$ws.Range("B9").Value = 0.008245133366701168
$ws.Range("D9").Value = 0.006739666653612321
$ws.Range("E9").Value = 0.02721138937451249
$ws.Range("B10").Value = 0.04677943709955017
$ws.Range("D10").Value = 0.000004853988798695153
$ws.Range("E10").Value = 0.04532007402961995
$ws.Range("B11").Value = 0.01116923574030743
$ws.Range("D11").Value = 0.008478677509239843
$ws.Range("E11").Value = 0.02881084924443909
$ws.Range("B12").Value = 0.05117677560327184
$ws.Range("D12").Value = 0.000005386936448215256
$ws.Range("E12").Value = 0.04401841069799114
$ws.Range("B13").Value = 0.07539307529505576
$ws.Range("D13").Value = 0.000004527619054415351
$ws.Range("E13").Value = 0.05040234842216126
$ws.Range("B14").Value = 0.0541008779768781
$ws.Range("D14").Value = 0.000005740966499468505
$ws.Range("E14").Value = 0.04335567385030175
$ws.Range("B15").Value = 0.07979041379877744
$ws.Range("D15").Value = 0.000004847413451221902
$ws.Range("E15").Value = 0.04877508126484811
$ws.Range("B16").Value = 0.00008947954649831365
$ws.Range("D16").Value = 0.07488283394358181
$ws.Range("E16").Value = 0.06106354281609668
$ws.Range("B17").Value = 0.007410920423826244
$ws.Range("D17").Value = 0.4357284991145118
$ws.Range("E17").Value = 0.0373952052171247
$ws.Range("B18").Value = 0.04741846028679066
$ws.Range("D18").Value = 0.00000458813943862233
$ws.Range("E18").Value = 0.05748488084646997
$ws.Range("B19").Value = 0.07163475997857459
$ws.Range("D19").Value = 0.000004048237816401805
$ws.Range("E19").Value = 0.06468671433290586
$ws.Range("B20").Value = 0.05034256266039692
$ws.Range("D20").Value = 0.000004942240423239806
$ws.Range("E20").Value = 0.05446732071311675
$ws.Range("B21").Value = 0.07603209848229626
$ws.Range("D21").Value = 0.000004368070663818072
$ws.Range("E21").Value = 0.05985487191177444
$ws.Range("B22").Value = 0.01473236130115418
$ws.Range("D22").Value = 0.4543452846015076
$ws.Range("E22").Value = 0.03727236173491121
$ws.Range("B23").Value = 0.07895620085590252
$ws.Range("D23").Value = 0.000004580531530054096
$ws.Range("E23").Value = 0.05747673365531459
$ws.Range("B24").Value = 0.09244412290484241
$ws.Range("D24").Value = 0.000004043673032953418
$ws.Range("E24").Value = 0.06469343119460909

# --- Sheet: criteria ---
$ws = $wb.Worksheets.Item("criteria")
$ws.Range("B2").Value = 0.2771882676036963
$ws.Range("B3").Value = 0.07956211810363849
$ws.Range("B4").Value = 0.07596432089007199
$ws.Range("B5").Value = 0.0887594780761196
$ws.Range("B6").Value = 0.02904115627477501
$ws.Range("B7").Value = 0.03630892285481715
$ws.Range("B8").Value = 0.1080909842885875
$ws.Range("B9").Value = 0.05279605036889801
$ws.Range("B10").Value = 0.252288701539396

# --- Sheet: result ---
$ws = $wb.Worksheets.Item("result")
$ws.Range("B2").Value = 0.1240851971693286
$ws.Range("B3").Value = 0.05909586467733745
$ws.Range("B4").Value = 0.05305810104732629
$ws.Range("A5").Value = 6
$ws.Range("B5").Value = 0.04885343389373953
$ws.Range("A6").Value = 21
$ws.Range("B6").Value = 0.0483336870525143
$ws.Range("A7").Value = 13
$ws.Range("B7").Value = 0.04747779785223569
$ws.Range("A8").Value = 11
$ws.Range("B8").Value = 0.04645421345764324
$ws.Range("A9").Value = 19
$ws.Range("B9").Value = 0.04604120066266992
$ws.Range("A10").Value = 17
$ws.Range("B10").Value = 0.04524206408060186
$ws.Range("A11").Value = 12
$ws.Range("B11").Value = 0.04033785525800364
$ws.Range("A12").Value = 5
$ws.Range("B12").Value = 0.03995929507450195
$ws.Range("A13").Value = 3
$ws.Range("B13").Value = 0.03915635796245084
$ws.Range("A14").Value = 18
$ws.Range("B14").Value = 0.03902585477222233
$ws.Range("A15").Value = 10
$ws.Range("B15").Value = 0.03789172916338462
$ws.Range("A16").Value = 8
$ws.Range("B16").Value = 0.03711911656684169
$ws.Range("A17").Value = 0
$ws.Range("B17").Value = 0.03698702234251269
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = 0.03680366320696893
$ws.Range("A19").Value = 4
$ws.Range("B19").Value = 0.03542152758696333
$ws.Range("A20").Value = 14
$ws.Range("B20").Value = 0.02989529345926219
$ws.Range("A21").Value = 9
$ws.Range("B21").Value = 0.02904117403357469
$ws.Range("B22").Value = 0.02806177221136825
$ws.Range("A23").Value = 1
$ws.Range("B24").Value = 0.02519197384996366
$ws.Range("B23").Value = 0.02646580541858433
